$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.15140533333333
$ws.Range("H2").Value = 30.454216
$ws.Range("I2").Value = 0.4279451754041443
$ws.Range("J2").Value = 0.4279451754041442
$ws.Range("M2").Value = 1.695786666666667
$ws.Range("N2").Value = 5.087359999999999
$ws.Range("O2").Value = 0.1708240824160653
$ws.Range("P2").Value = 0.1708240824160653
$ws.Range("Q2").Value = 17.21461781219555
$ws.Range("R2").Value = 154.93156030976
$ws.Range("S2").Value = 0.07310334191279506
$ws.Range("T2").Value = 0.07310334191279505
$ws.Range("G3").Value = 10.15140533333333
$ws.Range("H3").Value = 30.454216
$ws.Range("I3").Value = 0.4279451754041443
$ws.Range("J3").Value = 0.4279451754041442
$ws.Range("O3").Value = 0.3337244291137108
$ws.Range("P3").Value = 0.3337244291137108
$ws.Range("Q3").Value = 33.63072946467289
$ws.Range("R3").Value = 302.676565182056
$ws.Range("S3").Value = 0.1428157593537148
$ws.Range("T3").Value = 0.1428157593537148
$ws.Range("G4").Value = 10.15140533333333
$ws.Range("H4").Value = 30.454216
$ws.Range("I4").Value = 0.4279451754041443
$ws.Range("J4").Value = 0.4279451754041442
$ws.Range("M4").Value = 3.148025333333333
$ws.Range("N4").Value = 9.444075999999999
$ws.Range("O4").Value = 0.317114498869273
$ws.Range("P4").Value = 0.3171144988692729
$ws.Range("Q4").Value = 31.95688115826844
$ws.Range("R4").Value = 287.611930424416
$ws.Range("S4").Value = 0.1357076198418083
$ws.Range("T4").Value = 0.1357076198418083
$ws.Range("G5").Value = 10.15140533333333
$ws.Range("H5").Value = 30.454216
$ws.Range("I5").Value = 0.4279451754041443
$ws.Range("J5").Value = 0.4279451754041442
$ws.Range("M5").Value = 1.770368
$ws.Range("N5").Value = 5.311104
$ws.Range("O5").Value = 0.178336989600951
$ws.Range("P5").Value = 0.178336989600951
$ws.Range("Q5").Value = 17.97172315716267
$ws.Range("R5").Value = 161.745508414464
$ws.Range("S5").Value = 0.07631845429582605
$ws.Range("T5").Value = 0.07631845429582602
$ws.Range("G6").Value = 7.459653666666667
$ws.Range("I6").Value = 0.3144710207121242
$ws.Range("J6").Value = 0.3144710207121242
$ws.Range("M6").Value = 1.695786666666667
$ws.Range("N6").Value = 5.087359999999999
$ws.Range("O6").Value = 0.1708240824160653
$ws.Range("P6").Value = 0.1708240824160653
$ws.Range("Q6").Value = 12.64998122588444
$ws.Range("R6").Value = 113.84983103296
$ws.Range("S6").Value = 0.05371922355959208
$ws.Range("T6").Value = 0.05371922355959207
$ws.Range("G7").Value = 7.459653666666667
$ws.Range("I7").Value = 0.3144710207121242
$ws.Range("J7").Value = 0.3144710207121242
$ws.Range("O7").Value = 0.3337244291137108
$ws.Range("P7").Value = 0.3337244291137108
$ws.Range("Q7").Value = 24.71318858090011
$ws.Range("R7").Value = 222.418697228101
$ws.Range("S7").Value = 0.1049466618599596
$ws.Range("T7").Value = 0.1049466618599596
$ws.Range("G8").Value = 7.459653666666667
$ws.Range("I8").Value = 0.3144710207121242
$ws.Range("J8").Value = 0.3144710207121242
$ws.Range("M8").Value = 3.148025333333333
$ws.Range("N8").Value = 9.444075999999999
$ws.Range("O8").Value = 0.317114498869273
$ws.Range("P8").Value = 0.3171144988692729
$ws.Range("Q8").Value = 23.48317872055955
$ws.Range("R8").Value = 211.348608485036
$ws.Range("S8").Value = 0.09972332014203401
$ws.Range("T8").Value = 0.099723320142034
$ws.Range("G9").Value = 7.459653666666667
$ws.Range("I9").Value = 0.3144710207121242
$ws.Range("J9").Value = 0.3144710207121242
$ws.Range("M9").Value = 1.770368
$ws.Range("N9").Value = 5.311104
$ws.Range("O9").Value = 0.178336989600951
$ws.Range("P9").Value = 0.178336989600951
$ws.Range("Q9").Value = 13.20633214254934
$ws.Range("R9").Value = 118.856989282944
$ws.Range("S9").Value = 0.05608181515053855
$ws.Range("T9").Value = 0.05608181515053854
$ws.Range("G10").Value = 2.030652666666667
$ws.Range("H10").Value = 6.091958
$ws.Range("I10").Value = 0.08560470034312097
$ws.Range("J10").Value = 0.08560470034312095
$ws.Range("M10").Value = 1.695786666666667
$ws.Range("N10").Value = 5.087359999999999
$ws.Range("O10").Value = 0.1708240824160653
$ws.Range("P10").Value = 0.1708240824160653
$ws.Range("Q10").Value = 3.443553716764444
$ws.Range("R10").Value = 30.99198345088
$ws.Range("S10").Value = 0.01462334438661587
$ws.Range("T10").Value = 0.01462334438661587
$ws.Range("G11").Value = 2.030652666666667
$ws.Range("H11").Value = 6.091958
$ws.Range("I11").Value = 0.08560470034312097
$ws.Range("J11").Value = 0.08560470034312095
$ws.Range("O11").Value = 0.3337244291137108
$ws.Range("P11").Value = 0.3337244291137108
$ws.Range("Q11").Value = 6.727376971653111
$ws.Range("R11").Value = 60.546392744878
$ws.Range("S11").Value = 0.02856837975145832
$ws.Range("T11").Value = 0.02856837975145832
$ws.Range("G12").Value = 2.030652666666667
$ws.Range("H12").Value = 6.091958
$ws.Range("I12").Value = 0.08560470034312097
$ws.Range("J12").Value = 0.08560470034312095
$ws.Range("M12").Value = 3.148025333333333
$ws.Range("N12").Value = 9.444075999999999
$ws.Range("O12").Value = 0.317114498869273
$ws.Range("P12").Value = 0.3171144988692729
$ws.Range("Q12").Value = 6.392546037867555
$ws.Range("R12").Value = 57.53291434080799
$ws.Range("S12").Value = 0.02714649165016308
$ws.Range("T12").Value = 0.02714649165016307
$ws.Range("G13").Value = 2.030652666666667
$ws.Range("H13").Value = 6.091958
$ws.Range("I13").Value = 0.08560470034312097
$ws.Range("J13").Value = 0.08560470034312095
$ws.Range("M13").Value = 1.770368
$ws.Range("N13").Value = 5.311104
$ws.Range("O13").Value = 0.178336989600951
$ws.Range("P13").Value = 0.178336989600951
$ws.Range("Q13").Value = 3.595002500181334
$ws.Range("R13").Value = 32.355022501632
$ws.Range("S13").Value = 0.01526648455488369
$ws.Range("T13").Value = 0.01526648455488369
$ws.Range("G14").Value = 4.079563666666666
$ws.Range("H14").Value = 12.238691
$ws.Range("I14").Value = 0.1719791035406106
$ws.Range("J14").Value = 0.1719791035406106
$ws.Range("M14").Value = 1.695786666666667
$ws.Range("N14").Value = 5.087359999999999
$ws.Range("O14").Value = 0.1708240824160653
$ws.Range("P14").Value = 0.1708240824160653
$ws.Range("Q14").Value = 6.91806967175111
$ws.Range("R14").Value = 62.26262704575999
$ws.Range("S14").Value = 0.0293781725570623
$ws.Range("T14").Value = 0.0293781725570623
$ws.Range("G15").Value = 4.079563666666666
$ws.Range("H15").Value = 12.238691
$ws.Range("I15").Value = 0.1719791035406106
$ws.Range("J15").Value = 0.1719791035406106
$ws.Range("O15").Value = 0.3337244291137108
$ws.Range("P15").Value = 0.3337244291137108
$ws.Range("Q15").Value = 13.51524222533678
$ws.Range("R15").Value = 121.637180028031
$ws.Range("S15").Value = 0.05739362814857803
$ws.Range("T15").Value = 0.05739362814857803
$ws.Range("G16").Value = 4.079563666666666
$ws.Range("H16").Value = 12.238691
$ws.Range("I16").Value = 0.1719791035406106
$ws.Range("J16").Value = 0.1719791035406106
$ws.Range("M16").Value = 3.148025333333333
$ws.Range("N16").Value = 9.444075999999999
$ws.Range("O16").Value = 0.317114498869273
$ws.Range("P16").Value = 0.3171144988692729
$ws.Range("Q16").Value = 12.84256977161289
$ws.Range("R16").Value = 115.583127944516
$ws.Range("S16").Value = 0.05453706723526755
$ws.Range("T16").Value = 0.05453706723526754
$ws.Range("G17").Value = 4.079563666666666
$ws.Range("H17").Value = 12.238691
$ws.Range("I17").Value = 0.1719791035406106
$ws.Range("J17").Value = 0.1719791035406106
$ws.Range("M17").Value = 1.770368
$ws.Range("N17").Value = 5.311104
$ws.Range("O17").Value = 0.178336989600951
$ws.Range("P17").Value = 0.178336989600951
$ws.Range("Q17").Value = 7.222328969429333
$ws.Range("R17").Value = 65.000960724864
$ws.Range("S17").Value = 0.03067023559970276
$ws.Range("T17").Value = 0.03067023559970276
